$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New shared string / header: PRED VERD in I40 (pushes old Diff_Real / Diff_IDE to J40 / K40) ---
$ws.Range("K40").Value = "Diff_IDE"
$ws.Range("J40").Value = "Diff_Real"
$ws.Range("I40").Value = "PRED VERD"

# --- Row 41 ---
$ws.Range("I41").Value = 691651.5
$ws.Range("J41").Formula = "=(100*(G41/H41))"
$ws.Range("K41").Formula = "=(100*(F41/H41))"
$ws.Range("L41").Formula = "=(100*(F41/I41))"
$ws.Range("M41").Value = 13
$ws.Range("N41").Value = 19

# --- Row 42 ---
$ws.Range("I42").Value = 1628478.6
$ws.Range("J42").Formula = "=(100*(G42/H42))"
$ws.Range("K42").Formula = "=(100*(F42/H42))"
$ws.Range("L42").Formula = "=(100*(F42/I42))"
$ws.Range("M42").Value = 5
$ws.Range("N42").Value = 2

# --- Row 43 ---
$ws.Range("I43").Value = 217975.98
$ws.Range("J43").Formula = "=(100*(G43/H43))"
$ws.Range("K43").Formula = "=(100*(F43/H43))"
$ws.Range("L43").Formula = "=(100*(F43/I43))"
$ws.Range("M43").Value = 18
$ws.Range("N43").Value = 30

# --- Row 44 ---
$ws.Range("I44").Value = 713190.40000000002
$ws.Range("J44").Formula = "=(100*(G44/H44))"
$ws.Range("K44").Formula = "=(100*(F44/H44))"
$ws.Range("L44").Formula = "=(100*(F44/I44))"
$ws.Range("M44").Value = 1
$ws.Range("N44").Value = 3

# --- Row 45 ---
$ws.Range("I45").Value = 320622
$ws.Range("J45").Formula = "=(100*(G45/H45))"
$ws.Range("K45").Formula = "=(100*(F45/H45))"
$ws.Range("L45").Formula = "=(100*(F45/I45))"
$ws.Range("M45").Value = 17
$ws.Range("N45").Value = 16

# --- Row 46 : weighted sums over rows 41-45 ---
$ws.Range("M46").Formula = "=(0.2*(M41+M42+M43+M44+M45))"
$ws.Range("N46").Formula = "=(0.2*(N41+N42+N43+N44+N45))"

# --- Row 48 : weighted sums over rows 43-46 ---
$ws.Range("M48").Formula = "=(0.2*(M43+M44+M45+M46))"
$ws.Range("N48").Formula = "=(0.2*(N43+N44+N45+N46))"

# --- Column I width (closest value this engine's pixel-quantized ColumnWidth can reach to 10.54296875) ---
$ws.Columns.Item(9).ColumnWidth = 9.6

# --- View / selection state ---
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H40:H45").Select()

$wb.Save()
